$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Title
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(1).Range.Text = "San Diego County's Gross Regional Product (GRP) Report"

# ---------------------------------------------------------------------------
# 2-4. The three narrative paragraphs that follow the title get replaced,
#      and a new paragraph (ERROR placeholder) is appended after them.
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(2).Range.Text = "To generate the San Diego County's Gross Regional Product (GRP) report, we will analyze the data from the provided CSV files: ``sec1-1_San Diego County_grp_data.csv`` and ``sec1-2_San Diego County_per_capita_grp_data.csv``. The analysis will focus on GRP trends from 2019 to 2023, both in total dollars and on a per capita basis."

$d.Paragraphs.Item(3).Range.Text = "From 2019 to 2023, San Diego County's total GRP increased from approximately `$244.28 billion to `$308.71 billion, representing a growth of about 26.4% over the five-year period. This upward trajectory is evident in the data, with the most significant year-over-year increase occurring between 2021 and 2022, when the GRP rose by approximately `$27.81 billion. Concurrently, the GRP per capita also showed a consistent upward trend, increasing from `$73,346.92 in 2019 to `$94,915.87 in 2023, a growth of approximately 29.4%. The largest increase in GRP per capita was observed between 2021 and 2022, with an increase of `$8,457.15."

$d.Paragraphs.Item(4).Range.Text = "ERROR: Chart data could not be parsed."

# New paragraph inserted right after paragraph 4, before the table.
$d.Paragraphs.Item(4).Range.InsertParagraphAfter()
$d.Paragraphs.Item(5).Range.Text = "These trends suggest a robust economic performance in San Diego County. Despite a slight decrease in population from 3,330,458 in 2019 to 3,252,468 in 2023, the county's economic output per resident has improved, indicating a strengthening economy. The consistent increase in GRP per capita implies that the economic benefits are being distributed effectively among the residents, driven by contributions from various industries. Overall, San Diego County's economy has demonstrated significant growth and resilience over the past five years."

# ---------------------------------------------------------------------------
# 5. Table: add a 4th ("Population") column, narrow every column from
#    2880 -> 2160 twips (108pt), and refresh header/data text.
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)
$t.Columns.Add() | Out-Null
for ($i = 1; $i -le 4; $i++) {
  $t.Columns.Item($i).Width = 108
}

$t.Cell(1,1).Range.Text = "Year"
$t.Cell(1,2).Range.Text = "Total GRP (in dollars)"
$t.Cell(1,3).Range.Text = "GRP Per Capita"
$t.Cell(1,4).Range.Text = "Population"

$tableRows = @(
  @("2019", "244,278,846,457.50", "73,346.92", "3,330,458"),
  @("2020", "244,822,303,264.04", "74,277.60", "3,300,000"),
  @("2021", "268,873,550,391.04", "82,099.95", "3,270,000"),
  @("2022", "296,683,894,891.99", "90,557.10", "3,260,000"),
  @("2023", "308,710,843,090.86", "94,915.87", "3,252,468")
)
for ($r = 0; $r -lt $tableRows.Count; $r++) {
  $rowIndex = $r + 2
  for ($c = 0; $c -lt 4; $c++) {
    $t.Cell($rowIndex, $c + 1).Range.Text = $tableRows[$r][$c]
  }
}

# ---------------------------------------------------------------------------
# 6. Caption, and removal of the trailing "In conclusion" paragraph.
# ---------------------------------------------------------------------------
$captionPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $para = $d.Paragraphs.Item($i)
  if ($para.Style.NameLocal -eq "Caption") {
    $captionPara = $para
  }
}
$captionPara.Range.Text = "San Diego County's GRP Data from 2019 to 2023"

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.Delete()

# ---------------------------------------------------------------------------
# 7. Footer & header text (use a sub-range Find so only the existing
#    default header/footer part is touched, and literal apostrophes are
#    preserved rather than being curled by Find's replacement pipeline).
# ---------------------------------------------------------------------------
$sec = $d.Sections.Item(1)

$ftr = $sec.Footers.Item(1)
$ftrFound = $ftr.Range.Duplicate
$ftrFound.Find.Execute("Page", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$ftrFound.Text = "Prepared by [Your Name]"

$hdr = $sec.Headers.Item(1)
$hdrFound = $hdr.Range.Duplicate
$hdrFound.Find.Execute("San Diego County GRP Analysis", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$hdrFound.Text = "San Diego County's GRP Report"

Write-Output "done"
